# Adds two new doctor rows (3 and 4) to the sheet.
# NumberFormat "@" is set before writing phone/mobile numbers so they are
# stored as text (matching the existing row 2 phone/mobile cells) instead
# of being auto-detected as numbers; Style is reset to "Normal" afterwards
# so no stray number-format style lingers on the cell. Website/Mobile
# cells that are blank in the source data ("" in the diff) are written as
# a lone text apostrophe, which yields an empty but still text-typed
# cell (matching the blank inlineStr cells in the target), then restyled
# to "Normal".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "ΠΑΠΑΔΑΚΗΣ ΠΑΥΛΟΣ"
$ws.Range("B3").Value = "Θεσσαλονίκης 160, Πλησίον Ηλεκτρικού Σταθμού, Αθήνα - Κάτω Πετράλωνα, 11853, ΑΤΤΙΚΗΣ"
$ws.Range("C3").Value = "Κλινική Δερματολογία - Αφροδισιολογία"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2103455493"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6932351230"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "pavderm@gmail.com"

$ws.Range("A4").Value = "ΠΙΝΙΩΤΗ Γ. ΠΑΝΑΓΙΩΤΑ MD"
$ws.Range("B4").Value = "Αγ. Τρύφωνος 22, Γλυφάδα, 16561, ΑΤΤΙΚΗΣ"
$ws.Range("C4").Value = "Άγχος – Φοβίες - Κατάθλιψη – Κρίσεις Πανικού - Διαταραχές Μνήμης - Διαταραχές Ύπνου - Άνοια - Ψυχωσική Συνδρομή - Ψυχώσεις -  Πένθος  - Διαταραχές Διάθεσης - Συμβουλευτική - Ψυχοθεραπεία - Πιστοποιητικά Ψυχικής Υγείας - Ψυχολογική Υποστήριξη"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6976973323"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "https://drpinioti.gr/"
$ws.Range("G4").Value = "panagiotapinioti@gmail.com"
